$wb = $excel.ActiveWorkbook

# --- 1. Clear the stray empty cell at B2 on "ODI Batting" ---
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B2").Value = ""

# --- 2. Insert the new "ODI Batting Extra" sheet right after "ODI Batting" ---
$ws = $wb.Worksheets.Add($null, $odiBatting)
$ws.Name = "ODI Batting Extra"

# Copy the header formatting (bold + border) from the existing header row
# so the new header row matches the look of the other sheets.
$odiBatting.Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

# Keep every value as text (matches the source data, which stores numbers
# like "4586" / "0" as text strings, not numeric values) -- except column B
# (BATTING_POSITION), whose two non-blank entries are stored as real numbers.
$ws.Range("A1:A5").NumberFormat = "@"
$ws.Range("B1").NumberFormat = "@"
$ws.Range("C1:F5").NumberFormat = "@"

# --- Header row ---
$ws.Cells.Item(1,1).Value = "MATCH_CODE"
$ws.Cells.Item(1,2).Value = "BATTING_POSITION"
$ws.Cells.Item(1,3).Value = "NUM_4"
$ws.Cells.Item(1,4).Value = "NUM_6"
$ws.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Cells.Item(1,6).Value = "MAN_OF_MATCH"

# Blank data cells are still present in the source as explicit empty-text
# cells (not simply absent), so write a bare "'" (empty text marker) into
# each one instead of leaving the cell untouched.

# --- Row 2 : match 4586 ---
$ws.Cells.Item(2,1).Value = "4586"
$ws.Cells.Item(2,2).Value = "'"
$ws.Cells.Item(2,3).Value = "'"
$ws.Cells.Item(2,4).Value = "'"
$ws.Cells.Item(2,5).Value = "'"
$ws.Cells.Item(2,6).Value = "NO"

# --- Row 3 : match 4590 ---
$ws.Cells.Item(3,1).Value = "4590"
$ws.Cells.Item(3,2).Value = "'"
$ws.Cells.Item(3,3).Value = "'"
$ws.Cells.Item(3,4).Value = "'"
$ws.Cells.Item(3,5).Value = "'"
$ws.Cells.Item(3,6).Value = "NO"

# --- Row 4 : match 4592 ---
$ws.Cells.Item(4,1).Value = "4592"
$ws.Cells.Item(4,2).Value = 5
$ws.Cells.Item(4,3).Value = "0"
$ws.Cells.Item(4,4).Value = "0"
$ws.Cells.Item(4,5).Value = "'"
$ws.Cells.Item(4,6).Value = "NO"

# --- Row 5 : match 4641 ---
$ws.Cells.Item(5,1).Value = "4641"
$ws.Cells.Item(5,2).Value = 6
$ws.Cells.Item(5,3).Value = "0"
$ws.Cells.Item(5,4).Value = "0"
$ws.Cells.Item(5,5).Value = "1.94%"
$ws.Cells.Item(5,6).Value = "NO"
